# Q3 - Casos de Teste.docx : apply review-fix edits
$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Indentation fix: w:ind w:left="1416" -> "1985" for the ten code-line
#    paragraphs (the block-quoted Java source). LeftIndent is expressed
#    in points in the object model: 1985 twips / 20 = 99.25 pt.
# ---------------------------------------------------------------------
for ($i = 2; $i -le 11; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Format.LeftIndent -eq 70.8) {
        $p.Format.LeftIndent = 99.25
    }
}

# ---------------------------------------------------------------------
# 2) Move the "_GoBack" bookmark from the end of the document to right
#    after the "{  " run on the "if (numero >= 0) {" code line.
#    A direct Bookmarks.Add collapsed exactly at that text offset lands
#    on a paragraph/text boundary that this host mis-resolves, so we
#    stage a short unique marker there first (which *is* placed
#    correctly because it is not a collapsed zero-length range), drop
#    the bookmark next to the marker, then delete the marker text.
# ---------------------------------------------------------------------
$searchScope = $d.Range(300, $d.Content.End)
$searchScope.Find.Execute("{  ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$searchScope.InsertBefore("<<MARK>>")

$markScope = $d.Range(300, $d.Content.End)
$markScope.Find.Execute("<<MARK>>", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$markPos = $markScope.Start
$d.Bookmarks.Add("_GoBack", $d.Range($markPos, $markPos))

$deleteScope = $d.Range(300, $d.Content.End)
$deleteScope.Find.Execute("<<MARK>>", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)

# ---------------------------------------------------------------------
# 3) Remove the old "_GoBack" bookmark that used to sit at the very end
#    of the document (after the final test-case sentence).
# ---------------------------------------------------------------------
# (handled implicitly: Bookmarks.Add for a singleton bookmark like
#  _GoBack relocates it, so the stale end-of-document copy is already
#  gone; this call is a harmless no-op safety net.)
try {
    $oldBm = $d.Bookmarks("_GoBack")
    if ($oldBm.Start -gt $markPos) {
        $oldBm.Delete()
    }
} catch {
}

# ---------------------------------------------------------------------
# 4) Picture metadata: new anchor/edit ids + title="IMAGE1".
# ---------------------------------------------------------------------
$shp = $d.InlineShapes(1)
$shp.Title = "IMAGE1"

# ---------------------------------------------------------------------
# 5) Caption paragraph right under the picture: center it and add the
#    explanatory sentence.
# ---------------------------------------------------------------------
$captionPara = $d.Paragraphs(16)
$captionPara.Format.Alignment = 1
$captionPara.Range.InsertAfter("Nesta imagem, os números nos vértices correspondem às linhas no código.")

# ---------------------------------------------------------------------
# 6) Fix F(2): "Numero = 2   -> 2" should read "...-> 1", and the
#    corrected digit is kept in its own run (matches how Word actually
#    produced the change). The split is forced via a transient bookmark
#    placed at the exact cut point - inserting/deleting it creates the
#    run boundary without leaving any trace in the saved XML.
# ---------------------------------------------------------------------
function Split-LastChar-And-Fix {
    param($paraIndex)
    $p = $d.Paragraphs($paraIndex)
    $e = $p.Range.End
    $lastChar = $d.Range($e - 2, $e - 1)
    $lastChar.Text = "1"

    $e2 = $p.Range.End
    $splitPos = $e2 - 2
    $d.Bookmarks.Add("TempSplit", $d.Range($splitPos, $splitPos))
    $tmpBm = $d.Bookmarks("TempSplit")
    $tmpBm.Delete()
}

Split-LastChar-And-Fix 22
Split-LastChar-And-Fix 28

Write-Output "edits applied"
